# Refresh the crypto price / Volume(1h) table with the latest pull.
# Also: rows 25/26 and 29/30/31 swapped rank position (coin name, link,
# price and volume all move together with their row).
# (commit: "Updated cryptos list ... with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.311.58'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '2.364.10'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'520.79"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = "'135.93"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = "'0.997"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  -0.21%  '
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('D14').Value = '2.788.73'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Value = '57.312.18'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '2.369.05'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = "'330.28"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').Value = "'6.72"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').Value = "'0.997"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').Value = "'9.06"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +15.98%  '
$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').Value = "'1.03"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.165"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('E27').Value = '  +12.48%  '
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = "'166.71"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.20%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = "'6.30"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'1.69"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  +3.28%  '
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').Value = "'0.920"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.90%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = "'1.61"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.65%  '
$ws.Range('D39').Value = "'38.75"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.98%  '
$ws.Range('D40').Value = "'150.13"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.39%  '
$ws.Range('D41').Value = "'0.386"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = "'292.49"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.93%  '
$ws.Range('D43').Value = "'3.66"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E44').Value = '  +2.83%  '
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('D46').Value = "'0.0510"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').Value = "'0.567"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = "'18.23"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.04%  '
$ws.Range('D49').Value = "'0.0220"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('D50').Value = "'17.73"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('E51').Value = '  +1.41%  '
